# Update the "F" column (view/visit counts) values on the sheets that hold
# the exhibition data: "展览" and "全部类型" (they carry duplicate data).
#
# Row -> (old value, new value) for column F
#   2  : 1832 -> 1833
#   7  : 1505 -> 1510
#   9  : 592  -> 593
#   11 : 97   -> 98
#   19 : 3552 -> 3563
#   20 : 424  -> 425
#   21 : 317  -> 318
#   23 : 113  -> 120
#   24 : 336  -> 337
#   26 : 1308 -> 1320

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1833
    7  = 1510
    9  = 593
    11 = 98
    19 = 3563
    20 = 425
    21 = 318
    23 = 120
    24 = 337
    26 = 1320
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
